$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1667508
$ws.Range("J17").Value = 1724991
$ws.Range("L17").Value = 5174973
$ws.Range("N17").Value = -5175309

$ws.Range("H33").Value = 6526.0625
$ws.Range("I33").Value = 11425.777
$ws.Range("J33").Value = 226.42857
$ws.Range("K33").Value = 11425.777
$ws.Range("L33").Value = 226.42857
$ws.Range("M33").Value = -11196.777
$ws.Range("N33").Value = -684.42857

$ws.Range("H70").Value = 3263.25
$ws.Range("J70").Value = 4481.923
$ws.Range("L70").Value = 13445.769
$ws.Range("N70").Value = -13985.769

$ws.Range("H73").Value = 3263.25
$ws.Range("J73").Value = 4481.923
$ws.Range("L73").Value = 13445.769
$ws.Range("N73").Value = -15317.769

$ws.Range("H94").Value = 5856.8823
$ws.Range("I94").Value = 3070.875
$ws.Range("J94").Value = 8333.333000000001
$ws.Range("K94").Value = 3070.875
$ws.Range("L94").Value = 8333.333000000001
$ws.Range("M94").Value = -2619.875
$ws.Range("N94").Value = -9235.333000000001

$ws.Range("H113").Value = 6656.1514
$ws.Range("I113").Value = 2645
$ws.Range("J113").Value = 9998.777
$ws.Range("K113").Value = 2645
$ws.Range("L113").Value = 9998.777
$ws.Range("M113").Value = 609
$ws.Range("N113").Value = -16506.777

$ws.Range("H116").Value = 2741.5
$ws.Range("I116").Value = 2612.25
$ws.Range("J116").Value = 3000
$ws.Range("K116").Value = 2612.25
$ws.Range("L116").Value = 3000
$ws.Range("M116").Value = 829.75
$ws.Range("N116").Value = -9884

$ws.Range("H135").Value = 537.13043
$ws.Range("I135").Value = 537.13043
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 4834.173870000001
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -2299.173870000001
$ws.Range("N135").ClearContents()


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1298.76
$ws.Range("I2").Value = 1415.7894
$ws.Range("K2").Value = 1415.7894
$ws.Range("M2").Value = -1302.7894

$ws.Range("H14").Value = 4400
$ws.Range("I14").Value = 4400
$ws.Range("K14").Value = 4400
$ws.Range("M14").Value = -4225

$ws.Range("H32").Value = 4393.4
$ws.Range("I32").Value = 3626.739
$ws.Range("K32").Value = 3626.739
$ws.Range("M32").Value = -3339.739

$ws.Range("H61").Value = 2081.6667
$ws.Range("I61").Value = 2045.6173
$ws.Range("K61").Value = 2045.6173
$ws.Range("M61").Value = -1833.6173

$ws.Range("H116").Value = 1298.76
$ws.Range("I116").Value = 1415.7894
$ws.Range("K116").Value = 1415.7894
$ws.Range("M116").Value = 878.2106000000001

$ws.Range("H132").Value = 1737
$ws.Range("I132").Value = 1246.4445
$ws.Range("J132").Value = 2416.2307
$ws.Range("K132").Value = 3739.3335
$ws.Range("L132").Value = 7248.6921
$ws.Range("M132").Value = -1209.3335
$ws.Range("N132").Value = -12308.6921

$ws.Range("H136").Value = 2081.6667
$ws.Range("I136").Value = 2045.6173
$ws.Range("K136").Value = 6136.8519
$ws.Range("M136").Value = -3586.8519


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1298.76
$ws.Range("I3").Value = 1415.7894
$ws.Range("K3").Value = 1415.7894
$ws.Range("M3").Value = -1301.7894

$ws.Range("H99").Value = 1280
$ws.Range("I99").Value = 1112.5
$ws.Range("J99").Value = 1428.8889
$ws.Range("K99").Value = 1112.5
$ws.Range("L99").Value = 1428.8889
$ws.Range("M99").Value = 385.5
$ws.Range("N99").Value = -4424.8889

$ws.Range("H105").Value = 1502.6666
$ws.Range("I105").Value = 1104.5714
$ws.Range("J105").Value = 2060
$ws.Range("K105").Value = 1104.5714
$ws.Range("L105").Value = 2060
$ws.Range("M105").Value = 642.4286
$ws.Range("N105").Value = -5554


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("N21").ClearContents()

$ws.Range("H31").Value = 2189.767
$ws.Range("I31").Value = 1135.814
$ws.Range("J31").Value = 3700.4333
$ws.Range("K31").Value = 1135.814
$ws.Range("L31").Value = 3700.4333
$ws.Range("M31").Value = -840.8140000000001
$ws.Range("N31").Value = -4290.433300000001

$ws.Range("H34").Value = 2189.767
$ws.Range("I34").Value = 1135.814
$ws.Range("J34").Value = 3700.4333
$ws.Range("K34").Value = 1135.814
$ws.Range("L34").Value = 3700.4333
$ws.Range("M34").Value = -933.8140000000001
$ws.Range("N34").Value = -4104.433300000001


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 4738.5
$ws.Range("I82").Value = 2004.3334
$ws.Range("J82").Value = 6379
$ws.Range("K82").Value = 6013.0002
$ws.Range("L82").Value = 19137
$ws.Range("M82").Value = -5607.0002
$ws.Range("N82").Value = -19949

$ws.Range("H85").Value = 4738.5
$ws.Range("I85").Value = 2004.3334
$ws.Range("J85").Value = 6379
$ws.Range("K85").Value = 6013.0002
$ws.Range("L85").Value = 19137
$ws.Range("M85").Value = -4609.0002
$ws.Range("N85").Value = -21945

$ws.Range("H102").Value = 10946923
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 10946923
$ws.Range("K102").Value = 0
$ws.Range("N102").Value = -32845637
$ws.Range("L102").Value = 32840769
$ws.Range("M102").ClearContents()

$ws.Range("H127").Value = 829
$ws.Range("J127").Value = 829
$ws.Range("L127").Value = 2487
$ws.Range("N127").Value = -12407

$ws.Range("H131").Value = 3203.7144
$ws.Range("J131").Value = 3768.5576
$ws.Range("L131").Value = 11305.6728
$ws.Range("N131").Value = -21385.6728


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6947630.5
$ws.Range("I122").Value = 12347387
$ws.Range("J122").Value = 5085.7144
$ws.Range("K122").Value = 37042161
$ws.Range("L122").Value = 15257.1432
$ws.Range("M122").Value = -37039711
$ws.Range("N122").Value = -20157.1432

$ws.Range("H132").Value = 1371.7
$ws.Range("I132").Value = 1109.3334
$ws.Range("J132").Value = 2490.2104
$ws.Range("K132").Value = 3328.0002
$ws.Range("L132").Value = 7470.6312
$ws.Range("M132").Value = -798.0001999999999
$ws.Range("N132").Value = -12530.6312

$ws.Range("H136").Value = 6411374.5
$ws.Range("I136").Value = 1139.6666
$ws.Range("K136").Value = 3418.9998
$ws.Range("M136").Value = -868.9998000000001


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4117523.8
$ws.Range("J136").Value = 4726.5713
$ws.Range("L136").Value = 14179.7139
$ws.Range("N136").Value = -19279.7139

